$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update isotope names in column B (fix 48Ca->48Cd and 132Sn->132Cs mislabeling) ---
$ws.Range("B2").Value = "16O 1p1/2"
$ws.Range("B3").Value = "16O 1d5/2"
$ws.Range("B4").Value = "16O 2s1/2"
$ws.Range("B5").Value = "40Ca 1d5/2"
$ws.Range("B6").Value = "40Ca 2s1/2"
$ws.Range("B7").Value = "40Ca 1d3/2"
$ws.Range("B8").Value = "40Ca 1f7/2"
$ws.Range("B9").Value = "40Ca 2p3/2"
$ws.Range("B10").Value = "40Ca 2p1/2"
$ws.Range("B11").Value = "40Ca 1f5/2"
$ws.Range("B12").Value = "48Cd 1d5/2"
$ws.Range("B13").Value = "48Cd 2s1/2"
$ws.Range("B14").Value = "48Cd 1d3/2"
$ws.Range("B15").Value = "48Cd 1f7/2"
$ws.Range("B16").Value = "48Cd 2p3/2"
$ws.Range("B17").Value = "48Cd 2p1/2"
$ws.Range("B18").Value = "48Cd 1f5/2"
$ws.Range("B19").Value = "56Ni 1f7/2"
$ws.Range("B20").Value = "56Ni 2p3/2"
$ws.Range("B21").Value = "56Ni 1f5/2"
$ws.Range("B22").Value = "56Ni 2p1/2"
$ws.Range("B23").Value = "100Sn 2p1/2"
$ws.Range("B24").Value = "100Sn 1g9/2"
$ws.Range("B25").Value = "100Sn 2d5/2"
$ws.Range("B26").Value = "100Sn 1g7/2"
$ws.Range("B27").Value = "100Sn 3s1/2"
$ws.Range("B28").Value = "100Sn 1h11/2"
$ws.Range("B29").Value = "100Sn 2d3/2"
$ws.Range("B30").Value = "132Cs 1g7/2"
$ws.Range("B31").Value = "132Cs 2d5/2"
$ws.Range("B32").Value = "132Cs 3s1/2"
$ws.Range("B33").Value = "132Cs 1h11/2"
$ws.Range("B34").Value = "132Cs 2d3/2"
$ws.Range("B35").Value = "132Cs 2f7/2"
$ws.Range("B36").Value = "132Cs 3p3/2"
$ws.Range("B37").Value = "132Cs 1h9/2"
$ws.Range("B38").Value = "132Cs 2f5/2"
$ws.Range("B39").Value = "208Pb 1h9/2"
$ws.Range("B40").Value = "208Pb 2f7/2"
$ws.Range("B41").Value = "208Pb 1i13/2"
$ws.Range("B42").Value = "208Pb 3p3/2"
$ws.Range("B43").Value = "208Pb 2f5/2"
$ws.Range("B44").Value = "208Pb 3p1/2"
$ws.Range("B45").Value = "208Pb 2g9/2"
$ws.Range("B46").Value = "208Pb 1i11/2"
$ws.Range("B47").Value = "208Pb 1j15/2"
$ws.Range("B48").Value = "208Pb 3d5/2"
$ws.Range("B49").Value = "208Pb 4s1/2"
$ws.Range("B50").Value = "208Pb 2g7/2"
$ws.Range("B51").Value = "16O 1p1/2"
$ws.Range("B52").Value = "16O 1d5/2"
$ws.Range("B53").Value = "16O 2s1/2"
$ws.Range("B54").Value = "16O 1d3/2"
$ws.Range("B55").Value = "40Ca 1d5/2"
$ws.Range("B56").Value = "40Ca 2s1/2"
$ws.Range("B57").Value = "40Ca 1d3/2"
$ws.Range("B58").Value = "40Ca 1f7/2"
$ws.Range("B59").Value = "40Ca 2p3/2"
$ws.Range("B60").Value = "48Cd 1d5/2"
$ws.Range("B61").Value = "48Cd 1d3/2"
$ws.Range("B62").Value = "48Cd 2s1/2"
$ws.Range("B63").Value = "48Cd 1f7/2"
$ws.Range("B64").Value = "48Cd 2p3/2"
$ws.Range("B65").Value = "48Cd 2p1/2"
$ws.Range("B66").Value = "56Ni 1f7/2"
$ws.Range("B67").Value = "56Ni 2p3/2"
$ws.Range("B68").Value = "56Ni 1f5/2"
$ws.Range("B69").Value = "56Ni 2p1/2"
$ws.Range("B70").Value = "100Sn 1f5/2"
$ws.Range("B71").Value = "100Sn 2p3/2"
$ws.Range("B72").Value = "100Sn 2p1/2"
$ws.Range("B73").Value = "100Sn 1g9/2"
$ws.Range("B74").Value = "100Sn 1g7/2"
$ws.Range("B75").Value = "132Cs 2p1/2"
$ws.Range("B76").Value = "132Cs 1g9/2"
$ws.Range("B77").Value = "132Cs 1g7/2"
$ws.Range("B78").Value = "132Cs 2d5/2"
$ws.Range("B79").Value = "132Cs 2d3/2"
$ws.Range("B80").Value = "132Cs 1h11/2"
$ws.Range("B81").Value = "208Pb 1g7/2"
$ws.Range("B82").Value = "208Pb 2d5/2"
$ws.Range("B83").Value = "208Pb 1h11/2"
$ws.Range("B84").Value = "208Pb 2d3/2"
$ws.Range("B85").Value = "208Pb 3s1/2"
$ws.Range("B86").Value = "208Pb 1h9/2"
$ws.Range("B87").Value = "208Pb 2f7/2"
$ws.Range("B88").Value = "208Pb 1i13/2"
$ws.Range("B89").Value = "208Pb 1f5/2"
$ws.Range("B90").Value = "208Pb 3p3/2"

# --- Add new column M header "B_init" ---
$ws.Range("M1").Value = "B_init"
$ws.Range("M1").Style = $ws.Range("L1").Style

# --- Populate new column M (B_init) data ---
$ws.Range("M2").Value = -15.66
$ws.Range("M3").Value = -4.14
$ws.Range("M4").Value = -3.27
$ws.Range("M5").Value = -22.39
$ws.Range("M6").Value = -18.19
$ws.Range("M7").Value = -15.64
$ws.Range("M8").Value = -8.36
$ws.Range("M9").Value = -5.84
$ws.Range("M10").Value = -4.2
$ws.Range("M11").Value = -1.56
$ws.Range("M12").Value = -15.61
$ws.Range("M13").Value = -12.55
$ws.Range("M14").Value = -12.53
$ws.Range("M15").Value = -10
$ws.Range("M16").Value = -4.6
$ws.Range("M17").Value = -2.86
$ws.Range("M18").Value = -1.2
$ws.Range("M19").Value = -16.64
$ws.Range("M20").Value = -10.25
$ws.Range("M21").Value = -9.48
$ws.Range("M22").Value = -9.13
$ws.Range("M23").Value = -18.38
$ws.Range("M24").Value = -17.93
$ws.Range("M25").Value = -11.13
$ws.Range("M26").Value = -10.93
$ws.Range("M27").Value = -9.3
$ws.Range("M28").Value = -8.6
$ws.Range("M29").Value = -9.2
$ws.Range("M30").Value = -9.75
$ws.Range("M31").Value = -8.97
$ws.Range("M32").Value = -7.64
$ws.Range("M33").Value = -7.54
$ws.Range("M34").Value = -7.31
$ws.Range("M35").Value = -2.47
$ws.Range("M36").Value = 3.07944307104533
$ws.Range("M37").Value = -0.86
$ws.Range("M38").Value = 1.70678887664858
$ws.Range("M39").Value = -11.4
$ws.Range("M40").Value = -9.81
$ws.Range("M41").Value = -9.24
$ws.Range("M42").Value = -8.26
$ws.Range("M43").Value = -7.94
$ws.Range("M44").Value = -6.30466253792879
$ws.Range("M45").Value = -3.94
$ws.Range("M46").Value = -3.16
$ws.Range("M47").Value = -2.51
$ws.Range("M48").Value = 1.47128438332538
$ws.Range("M49").Value = 2.4716196788387
$ws.Range("M50").Value = -0.445616823906819
$ws.Range("M51").Value = -12.13
$ws.Range("M52").Value = -0.6
$ws.Range("M53").Value = -0.11
$ws.Range("M54").Value = 4.688
$ws.Range("M55").Value = -15.07
$ws.Range("M56").Value = -10.92
$ws.Range("M57").Value = -8.33
$ws.Range("M58").Value = -1.09
$ws.Range("M59").Value = 0.69
$ws.Range("M60").Value = -21.47
$ws.Range("M61").Value = -16.18
$ws.Range("M62").Value = -16.1
$ws.Range("M63").Value = -9.35
$ws.Range("M64").Value = -6.44
$ws.Range("M65").Value = -4.64
$ws.Range("M66").Value = -7.17
$ws.Range("M67").Value = -0.69
$ws.Range("M68").Value = 0.33
$ws.Range("M69").Value = 0.41
$ws.Range("M70").Value = -8.71
$ws.Range("M71").Value = -6.38
$ws.Range("M72").Value = -3.53
$ws.Range("M73").Value = -2.92
$ws.Range("M74").Value = 2.2031135908367
$ws.Range("M75").Value = -16.01
$ws.Range("M76").Value = -15.71
$ws.Range("M77").Value = -9.68
$ws.Range("M78").Value = -8.72
$ws.Range("M79").Value = -3.33777056868252
$ws.Range("M80").Value = -6.89
$ws.Range("M81").Value = -12
$ws.Range("M82").Value = -9.82
$ws.Range("M83").Value = -9.36
$ws.Range("M84").Value = -0.606436334495733
$ws.Range("M85").Value = -0.823026550062739
$ws.Range("M86").Value = -3.8
$ws.Range("M87").Value = -0.788046205296604
$ws.Range("M88").Value = -2.1
$ws.Range("M89").Value = 2.05124484053625
$ws.Range("M90").Value = 2.51801186433608

# --- Apply 0.00 number format + recomputed values for K (xmatch) cells that needed correction ---
$ws.Range("K36").NumberFormat = "0.00"
$ws.Range("K36").Value = 0.0564024168460056
$ws.Range("M36").NumberFormat = "0.00"
$ws.Range("K38").NumberFormat = "0.00"
$ws.Range("K38").Value = 0.0115998378585317
$ws.Range("M38").NumberFormat = "0.00"
$ws.Range("K44").NumberFormat = "0.00"
$ws.Range("K44").Value = 0.149788966724562
$ws.Range("M44").NumberFormat = "0.00"
$ws.Range("K48").NumberFormat = "0.00"
$ws.Range("K48").Value = 0.0402153227605435
$ws.Range("M48").NumberFormat = "0.00"
$ws.Range("K49").NumberFormat = "0.00"
$ws.Range("K49").Value = 0.0135315877292142
$ws.Range("M49").NumberFormat = "0.00"
$ws.Range("K50").NumberFormat = "0.00"
$ws.Range("K50").Value = 0.0148771405461984
$ws.Range("M50").NumberFormat = "0.00"
$ws.Range("K74").NumberFormat = "0.00"
$ws.Range("K74").Value = 0.00773583300609256
$ws.Range("M74").NumberFormat = "0.00"
$ws.Range("K79").NumberFormat = "0.00"
$ws.Range("K79").Value = 0.0584156989007292
$ws.Range("M79").NumberFormat = "0.00"
$ws.Range("K84").NumberFormat = "0.00"
$ws.Range("K84").Value = 0.041984151677683
$ws.Range("M84").NumberFormat = "0.00"
$ws.Range("K85").NumberFormat = "0.00"
$ws.Range("K85").Value = 0.0425255091646171
$ws.Range("M85").NumberFormat = "0.00"
$ws.Range("K87").NumberFormat = "0.00"
$ws.Range("K87").Value = 0.0186124770960984
$ws.Range("M87").NumberFormat = "0.00"
$ws.Range("K89").NumberFormat = "0.00"
$ws.Range("K89").Value = 0.0414784573715631
$ws.Range("M89").NumberFormat = "0.00"
$ws.Range("K90").NumberFormat = "0.00"
$ws.Range("K90").Value = 0.0467763884005684
$ws.Range("M90").NumberFormat = "0.00"

# --- Fix row 87: L87 had a stray misplaced value in N87; move it into L87 and clear N87 ---
$ws.Range("L87").Value = -419612941.28389442
$ws.Range("N87").ClearContents()

# --- Column width autofit for newly meaningful columns K and M ---
$ws.Columns("K").AutoFit()
$ws.Columns("M").AutoFit()

# --- Restore view state (selection) ---
$ws.Range("N80").Select()
